$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("A2:B539")
$keyRange = $ws.Range("A2")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

$ws.Columns.Item(1).ColumnWidth = 36.08984375

$ws.Range("E6").Select()
